# Update crypto price/volume figures per latest GitHub Actions scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '25.866.79'
$ws.Range('E2').Value = '  -0.35%  '
$ws.Range('D3').Value = '1.583.82'
$ws.Range('E3').Value = '  -2.05%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.13'
$ws.Range('E5').Value = '  -0.71%  '
$ws.Range('E6').Value = '  -0.20%  '
$ws.Range('E7').Value = '  -2.16%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.07'
$ws.Range('E10').Value = '  -0.47%  '
$ws.Range('E11').Value = '  -0.09%  '
$ws.Range('D12').Value = '1.804.72'
$ws.Range('E12').Value = '  -1.99%  '
$ws.Range('D13').Value = '1.582.04'
$ws.Range('E13').Value = '  -0.77%  '
$ws.Range('E14').Value = '  -2.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.504'
$ws.Range('E15').Value = '  -2.53%  '
$ws.Range('D16').Value = '25.873.86'
$ws.Range('E16').Value = '  -0.34%  '
$ws.Range('D17').Value = '0.0₃0724'
$ws.Range('E17').Value = '  -0.91%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '59.93'
$ws.Range('E18').Value = '  -2.63%  '
$ws.Range('E19').Value = '  -0.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '192.56'
$ws.Range('E20').Value = '  +0.77%  '
$ws.Range('E21').Value = '  -0.70%  '
$ws.Range('E22').Value = '  -0.65%  '
$ws.Range('E23').Value = '  -1.13%  '
$ws.Range('E24').Value = '  -0.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '141.16'
$ws.Range('E25').Value = '  -1.36%  '
$ws.Range('E26').Value = '  -0.21%  '
$ws.Range('E27').Value = '  -1.95%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.09'
$ws.Range('E28').Value = '  -0.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.43'
$ws.Range('E29').Value = '  -2.51%  '
$ws.Range('E30').Value = '  -4.92%  '
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.12'
$ws.Range('E32').Value = '  +0.37%  '
$ws.Range('E33').Value = '  -1.87%  '
$ws.Range('E34').Value = '  +0.73%  '
$ws.Range('E35').Value = '  -2.14%  '
$ws.Range('D36').Value = '1.097.18'
$ws.Range('E36').Value = '  -2.19%  '
$ws.Range('E37').Value = '  -0.32%  '
$ws.Range('E38').Value = '  -1.74%  '
$ws.Range('E39').Value = '  -0.51%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.501'
$ws.Range('E40').Value = '  -2.34%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.778'
$ws.Range('E41').Value = '  -4.85%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.798'
$ws.Range('E42').Value = '  +5.71%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '93.25'
$ws.Range('E43').Value = '  -4.06%  '
$ws.Range('E44').Value = '  +1.28%  '
$ws.Range('D45').Value = '1.718.66'
$ws.Range('E45').Value = '  -1.92%  '
$ws.Range('E46').Value = '  -1.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.51'
$ws.Range('E47').Value = '  +1.72%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '53.04'
$ws.Range('E48').Value = '  -1.32%  '
$ws.Range('E49').Value = '  -1.32%  '
$ws.Range('E50').Value = '  -0.82%  '
$ws.Range('E51').Value = '  -0.21%  '
